# Korrektur der Messwerte in Spalte I (Zeilen 8 und 9) auf Blatt "1. Stage"
# sowie Aktualisierung der Zellauswahl.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1. Stage")

# Korrigierte Messwerte
$ws.Range("I8").Value = 0.020759
$ws.Range("I9").Value = 0.177633

# Aktuelle Selektion im Tabellenblatt aktualisieren
$ws.Activate()
$ws.Range("I10").Select()
